# Exercise-2.docx edits
#
# wdFindWrap: wdFindContinue = 1
# wdReplaceAll = 2, wdReplaceNone = 0
# Range.Collapse: wdCollapseStart = 1, wdCollapseEnd = 0

$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $findText"
    }
}

# ---------------------------------------------------------------------------
# 1. "It wasn't that the flowers were awful" paragraph: rewrite the opening
#    clause.
# ---------------------------------------------------------------------------
Replace-Text "It wasn’t that the flowers were awful" "It wasn’t that there was something wrong with the flowers"

# ---------------------------------------------------------------------------
# 2. Wine-pouring paragraph: "pour out half the bottle" -> "pour nearly the
#    whole damn bottle out"
# ---------------------------------------------------------------------------
Replace-Text "she’d managed to pour out half the bottle within two days" "she’d managed to pour nearly the whole damn bottle out within two days"

# ---------------------------------------------------------------------------
# 3. Remove the old "_GoBack" bookmark that sat between " the" and
#    " greeting card" by re-writing across it (the bookmark is dropped by
#    the rewrite).
# ---------------------------------------------------------------------------
Replace-Text "got to sweeping up the remains of the greeting card" "got to sweeping up the remains of the greeting card"

# ---------------------------------------------------------------------------
# 4. Step-monster paragraph rewrite.
# ---------------------------------------------------------------------------
Replace-Text "go fuck herself instead of staying" "‘go to Hell’ herself instead of letting her stay"
Replace-Text "wasn’t 16 anymore and that kind of behavior wouldn’t cut it" "wasn’t 16 anymore and that mouthing off to the woman wouldn’t cut it"

# ---------------------------------------------------------------------------
# 5. The lone-space paragraph right after that becomes a writer's note.
# ---------------------------------------------------------------------------
$notePara = $d.Paragraphs(13)
$noteRng = $notePara.Range
$noteRng.Collapse(1)
$noteRng.InsertBefore("[That paragraph is in desperate need of rewriting – Lay off the alcohol talk here. Gist: Francesca hates her step-mom but now has to deal with her]")

# ---------------------------------------------------------------------------
# 6. Insert a new blank paragraph right after that note paragraph.
# ---------------------------------------------------------------------------
$d.Paragraphs(13).Range.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 7. Fill the middle one of the three trailing blank paragraphs (before the
#    closing quote paragraph) with the new "drunk twice" text, including
#    re-inserting the _GoBack bookmark inside it.
# ---------------------------------------------------------------------------
$drunkIdx = -1
for ($i = 2; $i -le ($d.Paragraphs.Count - 1); $i++) {
    $cur = $d.Paragraphs($i).Range.Text.Trim()
    $prev = $d.Paragraphs($i - 1).Range.Text.Trim()
    $next = $d.Paragraphs($i + 1).Range.Text.Trim()
    if ($cur -eq "" -and $prev -eq "" -and $next -eq "") {
        $afterNext = $d.Paragraphs($i + 2).Range.Text
        if ($afterNext -like '*You know when I moved out here*') {
            $drunkIdx = $i
        }
    }
}
if ($drunkIdx -eq -1) {
    throw "Could not locate the blank paragraph to fill with the drunk-twice text."
}
$beforeBookmark = "I’ve been drunk just twice in my life, and the "
$afterBookmark = ". The first time was at my matric-dance after party – I ended up dancing on top of the bar, it wasn’t my finest evening – The second time was last week. "

# Insert the full sentence first (so the paragraph actually has a run to
# anchor a bookmark on -- an empty, run-less paragraph cannot host one).
$drunkRng = $d.Paragraphs($drunkIdx).Range
$drunkRng.Collapse(1)
$drunkRng.InsertAfter($beforeBookmark + $afterBookmark)

# Re-locate the paragraph (indices/ranges are stable here since nothing was
# removed) and work out where the bookmark belongs: right between the two
# chunks of text we just inserted.
$filledPara = $d.Paragraphs($drunkIdx).Range
$bookmarkPos = $filledPara.Start + $beforeBookmark.Length
$bookmarkRng = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRng)
